$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.984.28'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").Value = '1.641.00'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'212.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = "'23.57"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.50%  '
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("D10").Value = "'0.0614"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("E11").Value = '  +2.26%  '
$ws.Range("D12").Value = '1.873.63'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").Value = '1.635.55'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("E15").Value = '  +3.98%  '
$ws.Range("D16").Value = "'65.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.08%  '
$ws.Range("D17").Value = '27.974.17'
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").Value = "'234.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = "'10.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").Value = "'4.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("E24").Value = '  -2.47%  '
$ws.Range("D25").Value = "'151.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("D26").Value = "'6.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.50%  '
$ws.Range("D27").Value = "'15.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = "'3.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '1.420.85'
$ws.Range("E34").Value = '  -3.56%  '
$ws.Range("E35").Value = '  +2.30%  '
$ws.Range("E36").Value = '  +1.28%  '
$ws.Range("E37").Value = '  +1.11%  '
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").Value = "'0.905"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.55%  '
$ws.Range("E41").Value = '  +0.71%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("E43").Value = '  +6.99%  '
$ws.Range("D44").Value = "'66.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.05%  '
$ws.Range("D45").Value = "'5.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.87%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = '1.782.11'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").Value = "'87.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = "'7.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.19%  '
